# Updated cryptos list on Sat Jul  1 06:55:57 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) for each coin row.
# Values that look like plain numbers are prefixed with a leading apostrophe
# so Excel stores them as literal text (matching the original inlineStr
# cells) instead of silently converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.418.74"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "1.918.27"
$ws.Range("E3").Value = "  +2.24%  "
$ws.Range("D4").Value = "'0.9987"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'241.14"
$ws.Range("E5").Value = "  +1.64%  "
$ws.Range("D6").Value = "'0.9992"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").Value = "'0.4693"
$ws.Range("E7").Value = "  -1.29%  "
$ws.Range("D8").Value = "'0.2851"
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("D9").Value = "'0.06835"
$ws.Range("E9").Value = "  +5.32%  "
$ws.Range("D10").Value = "'108.83"
$ws.Range("E10").Value = "  +14.29%  "
$ws.Range("D11").Value = "'18.24"
$ws.Range("E11").Value = "  -1.87%  "
$ws.Range("D12").Value = "1.907.48"
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("D13").Value = "'0.07635"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").Value = "'5.195"
$ws.Range("E14").Value = "  +3.10%  "
$ws.Range("D15").Value = "'0.6560"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("D16").Value = "'290.24"
$ws.Range("E16").Value = "  -2.48%  "
$ws.Range("D17").Value = "30.414.12"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("D18").Value = "'0.000007647"
$ws.Range("E18").Value = "  +2.00%  "
$ws.Range("D19").Value = "'0.9999"
$ws.Range("E19").Value = "  -0.08%  "
# Avalanche (row 20): only Volume(1h) changes, Price stays "12.94"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").Value = "2.140.64"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("D22").Value = "'0.9992"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").Value = "'5.234"
$ws.Range("E23").Value = "  +2.10%  "
$ws.Range("D24").Value = "'6.215"
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("D25").Value = "'21.74"
$ws.Range("E25").Value = "  +11.15%  "
$ws.Range("D26").Value = "'167.69"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").Value = "'9.270"
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("D28").Value = "'2.041"
$ws.Range("E28").Value = "  +5.20%  "
$ws.Range("D29").Value = "'0.1071"
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("D30").Value = "'1.370"
$ws.Range("E30").Value = "  +1.54%  "
$ws.Range("D31").Value = "'4.148"
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").Value = "'3.952"
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("D33").Value = "'0.05035"
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("D34").Value = "'0.7387"
$ws.Range("E34").Value = "  +3.12%  "
$ws.Range("D35").Value = "'1.147"
$ws.Range("E35").Value = "  -1.47%  "
$ws.Range("D36").Value = "'2.745"
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("D37").Value = "'0.02040"
$ws.Range("E37").Value = "  +7.40%  "
$ws.Range("D38").Value = "'2.685"
$ws.Range("E38").Value = "  -0.56%  "
$ws.Range("D39").Value = "'2.048"
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("D40").Value = "'0.8755"
$ws.Range("E40").Value = "  -2.13%  "
$ws.Range("D41").Value = "'108.49"
$ws.Range("E41").Value = "  +1.66%  "
$ws.Range("D42").Value = "'5.845"
$ws.Range("E42").Value = "  +5.08%  "
$ws.Range("D43").Value = "'0.9994"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").Value = "'0.4214"
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("D45").Value = "'51.74"
$ws.Range("E45").Value = "  +22.10%  "
$ws.Range("D46").Value = "'67.64"
$ws.Range("E46").Value = "  +3.80%  "
$ws.Range("D47").Value = "'7.187"
$ws.Range("E47").Value = "  -1.56%  "
$ws.Range("D48").Value = "'9.202"
$ws.Range("E48").Value = "  +3.02%  "

# Rows 49 and 50 swap ranking order: Elrond moves up to row 49,
# Algorand moves down to row 50 (each keeps its own refreshed price/volume).
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'34.77"
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.1207"
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("D51").Value = "'0.3885"
$ws.Range("E51").Value = "  +2.87%  "
